$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row fixes
$ws.Range("C1").Value = "Phần trăm lợi nhuận"
$ws.Range("A2").Value = "Số thứ tự"
$ws.Range("D2").Value = "Giá bìa "

# Existing rows (3-8) - book name swaps
$ws.Range("B3").Value = "Conan"
$ws.Range("B4").Value = "Năm mươi Sắc thái"
$ws.Range("B5").Value = "Cho tôi một vé đi tuổi thơ"
$ws.Range("B6").Value = "Chiến Thắng Con Quỷ Trong Bạn"
$ws.Range("B7").Value = "Đất Rừng Phương Nam"
$ws.Range("B8").Value = "Harry Potter Và Hòn Đá Phù Thuỷ"

# New rows 9-11 filled in
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Đắc Nhân Tâm"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = 50000

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Harry Potter Và Phòng Chứa Bí Mật"
$ws.Range("C10").Value = 50
$ws.Range("D10").Value = 100000

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Nguyên Tắc Vàng Của Napoleon Hill"
$ws.Range("C11").Value = 50
$ws.Range("D11").Value = 80000

# Update selection to reflect the new active range
$ws.Range("A12:D22").Select()
